$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.534.69"
$ws.Range("E2").Value = "  +1.51%  "

# Row 3
$ws.Range("D3").Value = "1.652.88"
$ws.Range("E3").Value = "  +2.57%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "301.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.74%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3568"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "

# Row 10
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.31%  "

# Row 11
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.227"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.005"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.31%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.23%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.426"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.473"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.68%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001208"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "

# Row 17
$ws.Range("D17").Value = "1.649.09"
$ws.Range("E17").Value = "  +1.97%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06981"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.87%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.780"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.64%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.72%  "

# Row 24
$ws.Range("D24").Value = "23.567.54"
$ws.Range("E24").Value = "  +1.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.492"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.45%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.939"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.22%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.215"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.30%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("B31").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C31").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D31").Value = "1.840.87"
$ws.Range("E31").Value = "  +2.62%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.109"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +9.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.177"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +10.06%  "

# Row 34
$ws.Range("E34").Value = "  +2.38%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.040"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02769"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.63%  "

# Row 37
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08755"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.28%  "

# Row 38
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.044"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.80%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2468"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.94%  "

# Row 40
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "13.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.21%  "

# Row 41
$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06930"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.22%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6957"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.328"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.83%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6468"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.27%  "

# Row 46
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.280"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.942"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.44%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07902"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.189"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
